$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1, copying format from E1 (bold/border/centered header style)
$ws.Cells.Item(1, 5).Copy()
$ws.Cells.Item(1, 6).PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Cells.Item(1, 6).Value = "time_taken"

$timestamps = @(
    "2021-10-05 13:40:40.242091",
    "2021-10-05 13:40:40.242102",
    "2021-10-05 13:40:40.242105",
    "2021-10-05 13:40:40.242108",
    "2021-10-05 13:40:40.242111",
    "2021-10-05 13:40:40.242113",
    "2021-10-05 13:40:40.242116",
    "2021-10-05 13:40:40.242119",
    "2021-10-05 13:40:40.242121",
    "2021-10-05 13:40:40.242124",
    "2021-10-05 13:40:40.242127",
    "2021-10-05 13:40:40.242129",
    "2021-10-05 13:40:40.242132",
    "2021-10-05 13:40:40.242134",
    "2021-10-05 13:40:40.242137",
    "2021-10-05 13:40:40.242139",
    "2021-10-05 13:40:40.242142",
    "2021-10-05 13:40:40.242145",
    "2021-10-05 13:40:40.242147",
    "2021-10-05 13:40:40.242150",
    "2021-10-05 13:40:40.242152",
    "2021-10-05 13:40:40.242155",
    "2021-10-05 13:40:40.242157",
    "2021-10-05 13:40:40.242160",
    "2021-10-05 13:40:40.242163",
    "2021-10-05 13:40:40.242166"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
